$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Description unknown, completed 06/15/2023 06:00:44 EDT, by WPJTOWN1.The search returned: 8 events."

$ws.Range("L2").Value = "Gross Weight"
$ws.Range("M2").Value = "Tare Weight"
$ws.Range("N2").Value = "Net Weight"
$ws.Range("O2").Value = "Car_no"

$ws.Range("A3").Value = "ITFX"
$ws.Range("B3").Value = 9725
$ws.Range("C3").Value = "JOHNSTOWN"
$ws.Range("D3").Value = "CO"
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1812
$ws.Range("H3").Value = "Placed Actual"
$ws.Range("J3").Value = "JOHNSTOWN"
$ws.Range("K3").Value = "CO"
$ws.Range("L3").Value = 202700
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 202700
$ws.Range("O3").Value = "ITFX9725"

$ws.Range("A4").Value = "ITFX"
$ws.Range("B4").Value = 9728
$ws.Range("C4").Value = "JOHNSTOWN"
$ws.Range("D4").Value = "CO"
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1812
$ws.Range("H4").Value = "Placed Actual"
$ws.Range("I4").ClearContents()
$ws.Range("J4").Value = "JOHNSTOWN"
$ws.Range("K4").Value = "CO"
$ws.Range("L4").Value = 202950
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 202950
$ws.Range("O4").Value = "ITFX9728"

$ws.Range("A5").Value = "MWCX"
$ws.Range("B5").Value = 102555
$ws.Range("C5").Value = "JOHNSTOWN"
$ws.Range("D5").Value = "CO"
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1811
$ws.Range("H5").Value = "Placed Actual"
$ws.Range("I5").ClearContents()
$ws.Range("J5").Value = "LOVELAND"
$ws.Range("K5").Value = "CO"
$ws.Range("L5").Value = 281100
$ws.Range("M5").Value = 73600
$ws.Range("N5").Value = 207500
$ws.Range("O5").Value = "MWCX102555"

$ws.Range("A6").Value = "MWCX"
$ws.Range("B6").Value = 102276
$ws.Range("C6").Value = "JOHNSTOWN"
$ws.Range("D6").Value = "CO"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 12
$ws.Range("G6").Value = 1304
$ws.Range("H6").Value = "Placed Actual"
$ws.Range("I6").ClearContents()
$ws.Range("J6").Value = "LOVELAND"
$ws.Range("K6").Value = "CO"
$ws.Range("L6").Value = 280350
$ws.Range("M6").Value = 78900
$ws.Range("N6").Value = 201450
$ws.Range("O6").Value = "MWCX102276"

$ws.Range("A7").Value = "MWCX"
$ws.Range("B7").Value = 102166
$ws.Range("C7").Value = "JOHNSTOWN"
$ws.Range("D7").Value = "CO"
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 12
$ws.Range("G7").Value = 1304
$ws.Range("H7").Value = "Placed Actual"
$ws.Range("J7").Value = "LOVELAND"
$ws.Range("K7").Value = "CO"
$ws.Range("L7").Value = 282400
$ws.Range("M7").Value = 82000
$ws.Range("N7").Value = 200400
$ws.Range("O7").Value = "MWCX102166"

$ws.Range("A8").Value = "MWCX"
$ws.Range("B8").Value = 102330
$ws.Range("C8").Value = "LOVELAND"
$ws.Range("D8").Value = "CO"
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 12
$ws.Range("G8").Value = 1045
$ws.Range("H8").Value = "Junction Received"
$ws.Range("I8").Value = "BNSF"
$ws.Range("J8").Value = "LOVELAND"
$ws.Range("K8").Value = "CO"
$ws.Range("L8").Value = 284850
$ws.Range("M8").Value = 79300
$ws.Range("N8").Value = 205550
$ws.Range("O8").Value = "MWCX102330"

$ws.Range("A9").Value = "MWCX"
$ws.Range("B9").Value = 102553
$ws.Range("C9").Value = "MEMPHIS"
$ws.Range("D9").Value = "TN"
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 14
$ws.Range("G9").Value = 1300
$ws.Range("H9").Value = "Junction Received"
$ws.Range("I9").Value = "NS"
$ws.Range("J9").Value = "LOVELAND"
$ws.Range("K9").Value = "CO"
$ws.Range("L9").Value = 281050
$ws.Range("M9").Value = 73400
$ws.Range("N9").Value = 207650
$ws.Range("O9").Value = "MWCX102553"

$ws.Range("A10").Value = "MWCX"
$ws.Range("B10").Value = 102328
$ws.Range("C10").Value = "VALDOSTA"
$ws.Range("D10").Value = "GA"
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 14
$ws.Range("G10").Value = 2248
$ws.Range("H10").Value = "Arrive In-Transit"
$ws.Range("J10").Value = "LOVELAND"
$ws.Range("K10").Value = "CO"
$ws.Range("L10").Value = 280550
$ws.Range("M10").Value = 79500
$ws.Range("N10").Value = 201050
$ws.Range("O10").Value = "MWCX102328"

$ws.AutoFilterMode = $false
foreach ($n in $wb.Names) {
    $n.Delete()
}

[void]$ws.Range("O3:O10").Select()
